$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '65.009.82'
$ws.Cells.Item(2, 5).Value = '  +2.30%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.151.61'
$ws.Cells.Item(3, 5).Value = '  +3.46%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.02%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''571.18'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +3.06%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''149.94'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +5.79%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.10%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '3.147.46'
$ws.Cells.Item(8, 5).Value = '  +3.35%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '''0.525'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +4.21%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''0.161'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +4.03%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''6.23'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +2.88%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '''0.506'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +6.43%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +16.94%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '''38.10'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +9.85%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '3.665.86'
$ws.Cells.Item(15, 5).Value = '  +3.57%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '65.069.26'
$ws.Cells.Item(16, 5).Value = '  +2.29%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '''7.21'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +7.29%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '3.151.35'
$ws.Cells.Item(18, 5).Value = '  +3.60%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +0.91%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''511.01'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +7.06%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''14.88'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +6.08%  '

# Row 22
$ws.Cells.Item(22, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(22, 4).Value = '''15.77'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +11.76%  '

# Row 23
$ws.Cells.Item(23, 2).Value = 'Polygon'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(23, 4).Value = '''0.737'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +9.22%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''7.83'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +3.98%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''84.71'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +4.52%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  -0.04%  '

# Row 27
$ws.Cells.Item(27, 2).Value = 'RenderToken'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(27, 4).Value = '''9.03'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +13.72%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'PancakeSwap'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(28, 4).Value = '''2.90'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +4.51%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '''2.21'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +8.87%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '''27.87'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +6.49%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '''1.24'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +8.85%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'Stacks'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(32, 4).Value = '''2.76'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +14.01%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(33, 4).Value = '''1.00'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +0.01%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '''6.31'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +12.25%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''6.61'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +7.29%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''55.80'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +1.87%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''477.87'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +9.74%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '''0.0873'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +8.17%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'dogwifhat'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(39, 4).Value = '''3.06'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +7.30%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'VeChain'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(40, 4).Value = '''0.0420'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +3.50%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '3.129.66'
$ws.Cells.Item(41, 5).Value = '  +5.30%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '''8.63'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +4.91%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''0.120'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +6.61%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''2.50'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +17.76%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '''0.290'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +12.00%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '''29.37'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +5.10%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '0.0₃0581'
$ws.Cells.Item(47, 5).Value = '  +13.33%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  -0.04%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'ThetaToken'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Cells.Item(49, 4).Value = '''2.35'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +12.97%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'Stellar'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(50, 4).Value = '''0.116'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +2.90%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '''123.23'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +5.48%  '
